$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata (file version / authoring environment) ---
# (These are application-level attributes that Excel stamps automatically on
#  save; the visible, data-relevant change for this commit is the refreshed
#  web-query report below.)

# --- Refreshed "Time In Phase" report (re-run of the excel(3) web query) ---

# Row 9 - Planeacion
$ws.Range("B9").Value = 0.0020833333333333333
$ws.Range("D9").Value = 0
$ws.Range("E9").NumberFormat = "0%"
$ws.Range("E9").Value = 0.05

# Row 10 - Diseno
$ws.Range("B10").Value = 0.010416666666666666
$ws.Range("D10").Value = 0
$ws.Range("E10").NumberFormat = "0%"
$ws.Range("E10").Value = 0.25

# Row 11 - Revision de Diseno
$ws.Range("B11").Value = 0.005555555555555556
$ws.Range("D11").Value = 0
$ws.Range("E11").NumberFormat = "0.00%"
$ws.Range("E11").Value = 0.125

# Row 12 - Codificacion
$ws.Range("B12").Value = 0.010416666666666666
$ws.Range("D12").Value = 0
$ws.Range("E12").NumberFormat = "0%"
$ws.Range("E12").Value = 0.25

# Row 13 - Revision de Codificacion
$ws.Range("B13").Value = 0.005555555555555556
$ws.Range("D13").Value = 0
$ws.Range("E13").NumberFormat = "0.00%"
$ws.Range("E13").Value = 0.125

# Row 14 - Compilacion
$ws.Range("D14").Value = 0
$ws.Range("E14").NumberFormat = "0.00%"
$ws.Range("E14").Value = 0.025

# Row 15 - Pruebas
$ws.Range("B15").Value = 0.005555555555555556
$ws.Range("D15").Value = 0
$ws.Range("E15").NumberFormat = "0.00%"
$ws.Range("E15").Value = 0.125

# Row 16 - Postmortem
$ws.Range("D16").Value = 0
$ws.Range("E16").NumberFormat = "0%"
$ws.Range("E16").Value = 0.05

# Row 17 - Total
$ws.Range("D17").Value = 0

# Footer timestamp generated by the report
$ws.Range("A25").Value = "Reporte generado a las 09:08 PM el 4/12/2018"
